$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.42"
$ws.Range("E2").Value = "'-2.36%"
$ws.Range("G2").Value = "'14"
$ws.Range("D3").Value = "'37.41"
$ws.Range("E3").Value = "'-5.41%"
$ws.Range("G3").Value = "'14"
$ws.Range("D4").Value = "'5.117"
$ws.Range("E4").Value = "'0.05%"
$ws.Range("G4").Value = "'14"
$ws.Range("D5").Value = "'0.07857"
$ws.Range("E5").Value = "'-4.22%"
$ws.Range("G5").Value = "'14"
$ws.Range("D6").Value = "'1.975"
$ws.Range("E6").Value = "'-0.03%"
$ws.Range("G6").Value = "'14"
$ws.Range("D7").Value = "'4.345"
$ws.Range("E7").Value = "'1.64%"
$ws.Range("G7").Value = "'14"
$ws.Range("D8").Value = "'8.221"
$ws.Range("E8").Value = "'-0.39%"
$ws.Range("G8").Value = "'14"
$ws.Range("D9").Value = "'3.098"
$ws.Range("E9").Value = "'-5.93%"
$ws.Range("G9").Value = "'14"
$ws.Range("D10").Value = "'0.9278"
$ws.Range("E10").Value = "'-0.51%"
$ws.Range("G10").Value = "'14"
$ws.Range("D11").Value = "'0.1294"
$ws.Range("E11").Value = "'-9.66%"
$ws.Range("G11").Value = "'14"
$ws.Range("D12").Value = "'0.1894"
$ws.Range("E12").Value = "'-4.43%"
$ws.Range("G12").Value = "'14"
$ws.Range("D13").Value = "'0.08748"
$ws.Range("E13").Value = "'-4.17%"
$ws.Range("G13").Value = "'14"
$ws.Range("D14").Value = "'0.03434"
$ws.Range("E14").Value = "'-3.32%"
$ws.Range("G14").Value = "'14"
$ws.Range("D15").Value = "'0.09754"
$ws.Range("E15").Value = "'-0.60%"
$ws.Range("G15").Value = "'14"
$ws.Range("D16").Value = "'0.001398"
$ws.Range("E16").Value = "'0.07%"
$ws.Range("G16").Value = "'14"
$ws.Range("D17").Value = "'0.005862"
$ws.Range("E17").Value = "'-6.85%"
$ws.Range("G17").Value = "'14"
$ws.Range("E18").Value = "'1,779.14%"
$ws.Range("G18").Value = "'14"
$ws.Range("D19").Value = "'3.590"
$ws.Range("E19").Value = "'-1.94%"
$ws.Range("G19").Value = "'14"
$ws.Range("E20").Value = "'-0.83%"
$ws.Range("G20").Value = "'14"
$ws.Range("D21").Value = "'0.1284"
$ws.Range("E21").Value = "'-1.44%"
$ws.Range("G21").Value = "'14"
$ws.Range("D22").Value = "'5.005"
$ws.Range("E22").Value = "'2.30%"
$ws.Range("G22").Value = "'14"
$ws.Range("D23").Value = "'0.2498"
$ws.Range("E23").Value = "'2.12%"
$ws.Range("G23").Value = "'14"
$ws.Range("D24").Value = "'0.04313"
$ws.Range("E24").Value = "'-0.19%"
$ws.Range("G24").Value = "'14"
$ws.Range("D25").Value = "'0.001223"
$ws.Range("E25").Value = "'-0.08%"
$ws.Range("G25").Value = "'14"
$ws.Range("D26").Value = "'0.004592"
$ws.Range("E26").Value = "'-4.16%"
$ws.Range("G26").Value = "'14"
$ws.Range("E27").Value = "'176.75%"
$ws.Range("G27").Value = "'14"
$ws.Range("G28").Value = "'14"
$ws.Range("G29").Value = "'14"
$ws.Range("G30").Value = "'14"
$ws.Range("G31").Value = "'14"
$ws.Range("G32").Value = "'14"
$ws.Range("G33").Value = "'14"
$ws.Range("G34").Value = "'14"
$ws.Range("G35").Value = "'14"
$ws.Range("G36").Value = "'14"
$ws.Range("G37").Value = "'14"
$ws.Range("G38").Value = "'14"
$ws.Range("E39").Value = "'1.88%"
$ws.Range("G39").Value = "'14"
$ws.Range("D40").Value = "'0.04993"
$ws.Range("E40").Value = "'-4.98%"
$ws.Range("G40").Value = "'14"
$ws.Range("D41").Value = "'0.007479"
$ws.Range("G41").Value = "'14"
$ws.Range("D42").Value = "'0.009813"
$ws.Range("E42").Value = "'0.82%"
$ws.Range("G42").Value = "'14"
$ws.Range("D43").Value = "'0.1352"
$ws.Range("E43").Value = "'-1.97%"
$ws.Range("G43").Value = "'14"
$ws.Range("D44").Value = "'0.002092"
$ws.Range("E44").Value = "'-2.10%"
$ws.Range("G44").Value = "'14"
$ws.Range("D45").Value = "'0.008004"
$ws.Range("E45").Value = "'-18.27%"
$ws.Range("G45").Value = "'14"
$ws.Range("D46").Value = "'0.00006507"
$ws.Range("E46").Value = "'1.95%"
$ws.Range("G46").Value = "'14"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.25%"
$ws.Range("G47").Value = "'14"
$ws.Range("D48").Value = "'0.003001"
$ws.Range("E48").Value = "'8.56%"
$ws.Range("G48").Value = "'14"
$ws.Range("D49").Value = "'0.001204"
$ws.Range("E49").Value = "'0.32%"
$ws.Range("G49").Value = "'14"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.25%"
$ws.Range("G50").Value = "'14"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.25%"
$ws.Range("G51").Value = "'14"
